$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as plain text (they are formatted numbers
# like "56.891.43" using "." as thousands separator, or trailing-zero
# decimals like "11.00" that must not be reinterpreted as numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.846.92"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.343.75"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.26"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.09"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.99"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.763.69"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.836.74"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.342.74"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.44"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "327.12"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.11"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  +4.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.70"
$ws.Range("E25").Value = "  +12.53%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.31"
$ws.Range("E27").Value = "  +8.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.55"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.48"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.70"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.07"
$ws.Range("E40").Value = "  +8.44%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "282.53"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.23"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.560"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.40"
$ws.Range("E48").Value = "  +7.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0216"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.29"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.00"
$ws.Range("E51").Value = "  +1.26%  "
